$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update the report title (shared string shown in the merged A1:M1
#    banner cell).
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Données COVID-19 Valais 22.06.2020"

# ---------------------------------------------------------------------
# 2) Revised daily inputs for the existing rows 106-114 (only the raw
#    input cells that actually changed; B/H/J/K are formulas and
#    recompute on their own).
# ---------------------------------------------------------------------
$ws.Range("G106").Value = 11
$ws.Range("I106").Value = 1

$ws.Range("F107").Value = 2
$ws.Range("G107").Value = 11

$ws.Range("F108").Value = 1
$ws.Range("G108").Value = 11

$ws.Range("E109").Value = 3
$ws.Range("F109").Value = 1
$ws.Range("G109").Value = 10
$ws.Range("I109").Value = 2

$ws.Range("E110").Value = 3
$ws.Range("F110").Value = 1
$ws.Range("G110").Value = 10

$ws.Range("E111").Value = 3
$ws.Range("F111").Value = 1
$ws.Range("G111").Value = 10

$ws.Range("E112").Value = 3
$ws.Range("F112").Value = 1
$ws.Range("G112").Value = 10

$ws.Range("C113").Value = 2
$ws.Range("D113").Value = 1
$ws.Range("E113").Value = 2
$ws.Range("F113").Value = 1
$ws.Range("G113").Value = 12

$ws.Range("C114").Value = 1
$ws.Range("E114").Value = 2
$ws.Range("F114").Value = 1
$ws.Range("G114").Value = 12

# ---------------------------------------------------------------------
# 3) Three more days (22.06.2020 ships data through 25.06.2020, i.e.
#    serials 44002-44004). Insert 3 blank rows right before the
#    current last row (115) so the old row 115 - with its bottom-
#    border "end of table" formatting and shared formulas - slides
#    down to become the new last row (118).
# ---------------------------------------------------------------------
$ws.Rows("115:117").Insert()

# --- formatting for the 3 freshly inserted rows (115-117): reproduce
#     row 114's "normal" look (thin left/right borders per column,
#     centered, no fill) ---------------------------------------------
function Set-NormalRowFormat($row) {
    $ws.Cells.Item($row, 1).Borders.Item(7).LineStyle = 1   # A left
    $ws.Cells.Item($row, 1).Borders.Item(10).LineStyle = 1  # A right

    $ws.Cells.Item($row, 2).Borders.Item(7).LineStyle = 1   # B left
    $ws.Cells.Item($row, 2).Interior.ColorIndex = -4142

    $ws.Cells.Item($row, 3).Borders.Item(7).LineStyle = 1   # C left

    $ws.Cells.Item($row, 4).Borders.Item(7).LineStyle = 1   # D left

    $ws.Cells.Item($row, 5).Borders.Item(7).LineStyle = 1   # E left
    $ws.Cells.Item($row, 5).Borders.Item(10).LineStyle = 1  # E right

    $ws.Cells.Item($row, 6).Borders.Item(10).LineStyle = 1  # F right

    $ws.Cells.Item($row, 8).Borders.Item(7).LineStyle = 1   # H left
    $ws.Cells.Item($row, 8).Borders.Item(10).LineStyle = 1  # H right
    $ws.Cells.Item($row, 8).Interior.ColorIndex = -4142

    $ws.Cells.Item($row, 9).Borders.Item(10).LineStyle = 1  # I right

    $ws.Cells.Item($row, 10).Borders.Item(10).LineStyle = 1 # J right
    $ws.Cells.Item($row, 10).Interior.ColorIndex = -4142

    $ws.Cells.Item($row, 11).Borders.Item(10).LineStyle = 1 # K right
    $ws.Cells.Item($row, 11).Interior.ColorIndex = -4142

    $ws.Cells.Item($row, 12).Borders.Item(7).LineStyle = 1  # L left

    $ws.Cells.Item($row, 13).Borders.Item(7).LineStyle = 1  # M left
    $ws.Cells.Item($row, 13).Borders.Item(10).LineStyle = 1 # M right
}

Set-NormalRowFormat 115
Set-NormalRowFormat 116
Set-NormalRowFormat 117

# --- values & formulas for row 115 (stays part of the original
#     shared-formula groups conceptually; write plain formulas, the
#     engine will store them as regular per-cell formulas) ----------
$ws.Range("A115").Value = 44001
$ws.Range("B115").Formula = "=B114+C115"
$ws.Range("C115").Value = 9
$ws.Range("D115").Value = 0
$ws.Range("E115").Value = 2
$ws.Range("F115").Value = 1
$ws.Range("G115").Value = 11
$ws.Range("H115").Formula = "=G115+E115"
$ws.Range("I115").Value = 1
$ws.Range("J115").Formula = "=J114+K115"
$ws.Range("K115").Formula = "=L115+M115"
$ws.Range("L115").NumberFormat = "General"
$ws.Range("L115").Value = 0
$ws.Range("L115").NumberFormat = "@"
$ws.Range("M115").NumberFormat = "General"
$ws.Range("M115").Value = 0
$ws.Range("M115").NumberFormat = "@"

# --- row 116 ----------------------------------------------------------
$ws.Range("A116").Value = 44002
$ws.Range("B116").Formula = "=B115+C116"
$ws.Range("C116").Value = 3
$ws.Range("D116").Value = 1
$ws.Range("E116").Value = 2
$ws.Range("F116").Value = 1
$ws.Range("G116").Value = 12
$ws.Range("H116").Formula = "=G116+E116"
$ws.Range("I116").Value = 0
$ws.Range("J116").Formula = "=J115+K116"
$ws.Range("K116").Formula = "=L116+M116"
$ws.Range("L116").NumberFormat = "General"
$ws.Range("L116").Value = 0
$ws.Range("L116").NumberFormat = "@"
$ws.Range("M116").NumberFormat = "General"
$ws.Range("M116").Value = 0
$ws.Range("M116").NumberFormat = "@"

# --- row 117 ----------------------------------------------------------
$ws.Range("A117").Value = 44003
$ws.Range("B117").Formula = "=B116+C117"
$ws.Range("C117").Value = 1
$ws.Range("D117").Value = 1
$ws.Range("E117").Value = 2
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 13
$ws.Range("H117").Formula = "=G117+E117"
$ws.Range("I117").Value = 0
$ws.Range("J117").Formula = "=J116+K117"
$ws.Range("K117").Formula = "=L117+M117"
$ws.Range("L117").NumberFormat = "General"
$ws.Range("L117").Value = 0
$ws.Range("L117").NumberFormat = "@"
$ws.Range("M117").NumberFormat = "General"
$ws.Range("M117").Value = 0
$ws.Range("M117").NumberFormat = "@"

# --- row 118: this is the old row 115 that Insert() slid down to the
#     bottom, so it already carries the correct bottom-border style
#     and L118/M118 are already 0 (unchanged) - only refresh the
#     inputs that differ, and fix up the 4 formulas that the insert
#     left pointing at row 114 instead of row 117. ------------------
$ws.Range("A118").Value = 44004
$ws.Range("B118").Formula = "=B117+C118"
$ws.Range("C118").Value = 1
$ws.Range("D118").Value = 0
$ws.Range("E118").Value = 2
$ws.Range("F118").Value = 1
$ws.Range("G118").Value = 13
$ws.Range("H118").Formula = "=G118+E118"
$ws.Range("I118").Value = 0
$ws.Range("J118").Formula = "=J117+K118"
$ws.Range("K118").Formula = "=L118+M118"

# ---------------------------------------------------------------------
# 4) View: scrolled down, A118 selected.
# ---------------------------------------------------------------------
$ws.Range("A118").Select()
$excel.ActiveWindow.ScrollRow = 100
